$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.809.27"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "3.529.38"
$ws.Range("E3").Value = "  +1.01%  "

$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "'604.64"
$ws.Range("E5").Value = "  -0.86%  "

$ws.Range("D6").Value = "'196.42"
$ws.Range("E6").Value = "  +5.33%  "

$ws.Range("D7").Value = "'0.630"
$ws.Range("E7").Value = "  +0.43%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").Value = "'0.201"
$ws.Range("E9").Value = "  -5.06%  "

$ws.Range("D10").Value = "'0.649"
$ws.Range("E10").Value = "  -0.58%  "

$ws.Range("D11").Value = "'53.68"
$ws.Range("E11").Value = "  +0.85%  "

$ws.Range("D12").Value = "'0.0000303"
$ws.Range("E12").Value = "  -1.07%  "

$ws.Range("D13").Value = "'9.51"
$ws.Range("E13").Value = "  -1.17%  "

$ws.Range("D14").Value = "4.087.44"
$ws.Range("E14").Value = "  +0.89%  "

$ws.Range("D15").Value = "'597.23"
$ws.Range("E15").Value = "  -1.35%  "

$ws.Range("D16").Value = "69.961.50"
$ws.Range("E16").Value = "  +0.72%  "

$ws.Range("D17").Value = "'19.11"
$ws.Range("E17").Value = "  +0.99%  "

$ws.Range("D18").Value = "'12.72"
$ws.Range("E18").Value = "  +0.29%  "

$ws.Range("D19").Value = "3.540.81"
$ws.Range("E19").Value = "  +1.68%  "

$ws.Range("E20").Value = "  +1.05%  "

$ws.Range("E21").Value = "  +0.61%  "

$ws.Range("D22").Value = "'18.56"
$ws.Range("E22").Value = "  +6.99%  "

$ws.Range("D23").Value = "'5.26"
$ws.Range("E23").Value = "  +4.86%  "

$ws.Range("D24").Value = "'102.02"
$ws.Range("E24").Value = "  -2.62%  "

$ws.Range("E25").Value = "  -0.43%  "

$ws.Range("D26").Value = "'3.16"
$ws.Range("E26").Value = "  +3.43%  "

$ws.Range("D27").Value = "'10.92"
$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("D28").Value = "'9.61"
$ws.Range("E28").Value = "  -3.00%  "

$ws.Range("D29").Value = "'33.47"
$ws.Range("E29").Value = "  -0.76%  "

$ws.Range("D30").Value = "'7.10"
$ws.Range("E30").Value = "  +1.24%  "

$ws.Range("E31").Value = "  +11.37%  "

$ws.Range("D32").Value = "'12.52"
$ws.Range("E32").Value = "  +0.06%  "

$ws.Range("E33").Value = "  -1.34%  "

$ws.Range("D34").Value = "'63.10"
$ws.Range("E34").Value = "  -0.38%  "

$ws.Range("D35").Value = "0.0₃0868"
$ws.Range("E35").Value = "  +12.73%  "

$ws.Range("D36").Value = "3.725.28"
$ws.Range("E36").Value = "  +4.41%  "

$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.11%  "

$ws.Range("D38").Value = "'3.06"
$ws.Range("E38").Value = "  -2.77%  "

$ws.Range("D39").Value = "'3.63"
$ws.Range("E39").Value = "  +1.29%  "

$ws.Range("E40").Value = "  -0.82%  "

$ws.Range("E41").Value = "  -0.21%  "

$ws.Range("D42").Value = "'488.51"
$ws.Range("E42").Value = "  -6.78%  "

$ws.Range("E43").Value = "  -2.93%  "

$ws.Range("E44").Value = "  -1.42%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.141"
$ws.Range("E45").Value = "  -3.37%  "

$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "'2.84"
$ws.Range("E46").Value = "  -4.90%  "

$ws.Range("D47").Value = "'3.30"
$ws.Range("E47").Value = "  -1.23%  "

$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  +0.18%  "

$ws.Range("D49").Value = "'8.54"
$ws.Range("E49").Value = "  -3.63%  "

$ws.Range("D50").Value = "'0.000254"
$ws.Range("E50").Value = "  +4.97%  "

$ws.Range("D51").Value = "'131.06"
$ws.Range("E51").Value = "  -0.07%  "
